# "QM working. Data getting into Orch queue."
# Populate the new Queue Management settings on the Settings sheet and point
# the dispatch workbook path at the relative Data\ location.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Fill in the new settings rows (9-16) in the order the author appears to
# have typed them - this is also the order new shared strings get interned.
$ws.Range("A9").Value = "DeletePriorRunNewQueueItems "
$ws.Range("B9").Value = $true

$ws.Range("A11").Value = "ApplicationUserNames"
$ws.Range("A12").Value = "NextJobProcessName"

$ws.Range("A13").Value = "RetryCount"
$ws.Range("B13").Value = 3

$ws.Range("A14").Value = "RetryIntervalSeconds "
$ws.Range("B14").Value = 5

$ws.Range("A15").Value = "StartJobApiJsonBody"

$ws.Range("A16").Value = "EnableFaultOnError"
$ws.Range("B16").Value = $false

$ws.Range("A10").Value = "QueueItemReferenceKey "

# Point ExcelFilePath at the relative Data\ path instead of the absolute one.
$ws.Range("B5").Value = "Data\DispatchExcel.xlsx"

$ws.Range("B10").Value = "Reference"

# Best-fit columns A and B now that their content has changed.
$ws.Columns("A:A").AutoFit()
$ws.Columns("B:B").AutoFit()

# Trim the two now-superfluous trailing blank rows at the bottom of the sheet.
$ws.Rows("996:997").Delete()

# Leave the selection where the author last clicked.
$null = $ws.Range("A17").Select()
